{"js": "// Office.js (Word JavaScript API) script.\n// Applies the diff:\n//   1. Removes the \"Meta description: ...\" paragraph (2nd paragraph, right\n//      after the H1 title).\n//   2. Before the final (italic \"Create a feature image...\" prompt)\n//      paragraph, inserts a new bold paragraph carrying the page title, and\n//      rewrites the final paragraph's text to the meta-description copy\n//      (keeping its italic run formatting).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Step 1: delete the \"Meta description\" paragraph -----------------\n// It is the paragraph right after the H1 title, starting with the bold\n// \"Meta description\" label.\nlet metaPara = null;\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"Meta description\") === 0) {\n    metaPara = p;\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// --- Step 2: replace the last paragraph with two paragraphs ----------\n// Re-load paragraphs since the collection changed after the delete above.\nconst paragraphs2 = body.paragraphs;\nparagraphs2.load(\"items\");\nawait context.sync();\n\nconst items = paragraphs2.items;\nconst lastPara = items[items.length - 1];\nconst wholeRange = lastPara.getRange(\"Whole\");\n\n// Flat-OPC OOXML fragment: the bold title paragraph followed by the\n// (still italic) description paragraph, replacing the old image-prompt\n// paragraph entirely.\nconst flatOpcXml = `<?xml version=\"1.0\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Boxing Arena Free - Review of Dreamtech's Slot Game</w:t></w:r></w:p>\n          <w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the exciting world of Dreamtech's Boxing Arena slot game with our review. Play for free and experience the thrills for yourself!</w:t></w:r></w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nwholeRange.insertOoxml(flatOpcXml, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the diff:\n#   1. Removes the \"Meta description: ...\" paragraph (2nd paragraph, right\n#      after the H1 title).\n#   2. Before the final (italic \"Create a feature image...\" prompt)\n#      paragraph, inserts a new bold paragraph carrying the page title, and\n#      rewrites the final paragraph's text to the meta-description copy\n#      (keeping its italic run formatting).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Meta description\" paragraph ----------------------\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Meta description\")) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# --- Step 2: replace the last paragraph with the title + description -----\n$count = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($count)\n$lastRange = $lastPara.Range\n$lastRange.Delete()\n\n$flatOpcXml = '<?xml version=\"1.0\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n      '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n          '<w:body>' +\n            '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Boxing Arena Free - Review of Dreamtech''s Slot Game</w:t></w:r></w:p>' +\n            '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Discover the exciting world of Dreamtech''s Boxing Arena slot game with our review. Play for free and experience the thrills for yourself!</w:t></w:r></w:p>' +\n          '</w:body>' +\n        '</w:document>' +\n      '</pkg:xmlData>' +\n    '</pkg:part>' +\n  '</pkg:package>'\n\n$lastRange.InsertXML($flatOpcXml)\n"}
